$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Country" column (A). Everything shifts one column left,
# so old Continent/TotalCases/.../Critical (B..J) become the new A..I.
$ws.Columns("A").Delete()

# Drop the now-trailing stat columns (old K..O -> now J..N) that are no
# longer part of the table ("Tot Cases/1M pop", "Deaths/1M pop",
# "TotalTests", "Tests/1M pop", "Population").
$ws.Range("J1:N1").EntireColumn.Delete()

# Refresh the data values to the latest snapshot.
# Row 2 - North America
$ws.Range("B2").Value = 126699963
$ws.Range("D2").Value = 1633266
$ws.Range("F2").Value = 122586057
$ws.Range("H2").Value = 2480640

# Row 3 - Asia
$ws.Range("B3").Value = 217430622
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = 1546531
$ws.Range("E3").ClearContents()
$ws.Range("F3").Value = 201143899
$ws.Range("G3").Value = 23504
$ws.Range("H3").Value = 14740192
$ws.Range("I3").Value = 15363

# Row 4 - Europe
$ws.Range("B4").Value = 249364068
$ws.Range("D4").Value = 2061022
$ws.Range("F4").Value = 245341050
$ws.Range("G4").Value = 5305
$ws.Range("H4").Value = 1961996

# Row 5 - South America
$ws.Range("B5").Value = 68711998
$ws.Range("F5").Value = 66471722
$ws.Range("G5").ClearContents()
$ws.Range("H5").Value = 884385

# Row 6 - Australia/Oceania
$ws.Range("B6").Value = 14365175
$ws.Range("H6").Value = 148581

# Row 7 - Africa
$ws.Range("B7").Value = 12823522
$ws.Range("F7").Value = 12085479
$ws.Range("G7").ClearContents()
$ws.Range("H7").Value = 479278
